$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 57884

# Row 3
$ws.Range("B3").Value = 57884

# Row 4
$ws.Range("B4").Value = 57884

# Row 5
$ws.Range("B5").Value = 57884

# Row 6
$ws.Range("B6").Value = 79243

# Row 7
$ws.Range("B7").Value = 79243

# Row 8
$ws.Range("B8").Value = 79243

# Row 9
$ws.Range("B9").Value = 79243

# Row 11
$ws.Range("A11").Value = 130961179
$ws.Range("B11").Value = 79862
$ws.Range("E11").Value = 6453
$ws.Range("F11").Value = "Vedskivlav"
$ws.Range("G11").Value = "Hertelidea botryosa"
$ws.Range("H11").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M11").ClearContents()
$ws.Range("AC11").Value = "Ringhack på stam i bakgrund"

# Row 12
$ws.Range("A12").Value = 130961218
$ws.Range("B12").Value = 57884
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "äldre spår"
$ws.Range("AC12").ClearContents()

# Row 13
$ws.Range("B13").Value = 79243

# Row 14
$ws.Range("B14").Value = 79243

# Row 15
$ws.Range("B15").Value = 79243

# Row 16
$ws.Range("B16").Value = 79243

# Row 17
$ws.Range("A17").Value = 130960378
$ws.Range("B17").Value = 57884
$ws.Range("E17").Value = 100109
$ws.Range("F17").Value = "Tretåig hackspett"
$ws.Range("G17").Value = "Picoides tridactylus"
$ws.Range("H17").Value = "(Linnaeus, 1758)"
$ws.Range("M17").Value = "äldre spår"
$ws.Range("Q17").Value = 446272
$ws.Range("R17").Value = 6759739

# Row 18
$ws.Range("A18").Value = 130960789
$ws.Range("B18").Value = 79243
$ws.Range("E18").Value = 6425
$ws.Range("F18").Value = "Garnlav"
$ws.Range("G18").Value = "Alectoria sarmentosa"
$ws.Range("H18").Value = "(Ach.) Ach."
$ws.Range("Q18").Value = 446284
$ws.Range("R18").Value = 6759886
$ws.Range("AC18").ClearContents()

# Row 19
$ws.Range("B19").Value = 79243

# Row 20
$ws.Range("A20").Value = 130961956
$ws.Range("B20").Value = 79862
$ws.Range("E20").Value = 6453
$ws.Range("F20").Value = "Vedskivlav"
$ws.Range("G20").Value = "Hertelidea botryosa"
$ws.Range("H20").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("M20").ClearContents()
$ws.Range("Q20").Value = 446084
$ws.Range("R20").Value = 6759981
$ws.Range("AC20").Value = "Miljöbilder"

# Row 21
$ws.Range("B21").Value = 79243

# Row 22
$ws.Range("A22").Value = 130963976
$ws.Range("B22").Value = 79243
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = "Garnlav"
$ws.Range("G22").Value = "Alectoria sarmentosa"
$ws.Range("H22").Value = "(Ach.) Ach."
$ws.Range("Q22").Value = 445929
$ws.Range("R22").Value = 6760099
$ws.Range("Z22").Value = "14:08"
$ws.Range("AB22").Value = "14:08"
$ws.Range("AC22").Value = "Miljöbild"

# Row 23
$ws.Range("A23").Value = 130962722
$ws.Range("B23").Value = 79862
$ws.Range("Q23").Value = 446008
$ws.Range("R23").Value = 6759948

# Row 24
$ws.Range("A24").Value = 130962640
$ws.Range("B24").Value = 79862
$ws.Range("E24").Value = 6453
$ws.Range("F24").Value = "Vedskivlav"
$ws.Range("G24").Value = "Hertelidea botryosa"
$ws.Range("H24").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q24").Value = 446038
$ws.Range("R24").Value = 6759945
$ws.Range("Z24").Value = "10:26"
$ws.Range("AB24").Value = "10:26"
$ws.Range("AC24").ClearContents()

# Row 25
$ws.Range("A25").Value = 130961746
$ws.Range("B25").Value = 57881
$ws.Range("E25").Value = 100049
$ws.Range("F25").Value = "Spillkråka"
$ws.Range("G25").Value = "Dryocopus martius"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
$ws.Range("M25").Value = "färska spår"
$ws.Range("Q25").Value = 446098
$ws.Range("R25").Value = 6760061

# Row 26
$ws.Range("A26").Value = 130962090
$ws.Range("B26").Value = 79243
$ws.Range("E26").Value = 6425
$ws.Range("F26").Value = "Garnlav"
$ws.Range("G26").Value = "Alectoria sarmentosa"
$ws.Range("H26").Value = "(Ach.) Ach."
$ws.Range("M26").ClearContents()
$ws.Range("Q26").Value = 446080
$ws.Range("R26").Value = 6759960

# Row 27
$ws.Range("A27").Value = 130961461
$ws.Range("B27").Value = 79243
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("M27").ClearContents()
$ws.Range("Q27").Value = 446088
$ws.Range("R27").Value = 6760088
$ws.Range("Z27").Value = "10:26"
$ws.Range("AB27").Value = "10:26"

# Row 28
$ws.Range("A28").Value = 130961750
$ws.Range("B28").Value = 79243
$ws.Range("Q28").Value = 446098
$ws.Range("R28").Value = 6760061
$ws.Range("AC28").Value = "Rikligt i en radie av ca 50 meter"

# Row 29
$ws.Range("A29").Value = 130963807
$ws.Range("B29").Value = 57881
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 445932
$ws.Range("R29").Value = 6760079
$ws.Range("Z29").Value = "14:08"
$ws.Range("AB29").Value = "14:08"
$ws.Range("AC29").ClearContents()

# Row 30
$ws.Range("B30").Value = 79833

# Row 31
$ws.Range("B31").Value = 79243

# Row 32
$ws.Range("B32").Value = 79243

# Row 33
$ws.Range("B33").Value = 79243

# Row 34
$ws.Range("B34").Value = 79243

# Row 35
$ws.Range("B35").Value = 57884
